# Fill in the two "future feature" slides (Blogs / Assignments) that were
# left with an empty content placeholder, with their bullet copy.
#
# Slide 12 - "Blogs" : Content Placeholder 2 (shape 1)
# Slide 13 - "Assignments" : Content Placeholder 2 (shape 1)

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 12 ("Blogs")
# ---------------------------------------------------------------------
$slide12 = $p.Slides.Item(12)
$body12 = $slide12.Shapes.Item("Content Placeholder 2").TextFrame.TextRange

$body12.Text = "Draft capabilities"
[void]$body12.InsertAfter("`rPost  management")
[void]$body12.InsertAfter("`rRich text ")
[void]$body12.InsertAfter("e")
[void]$body12.InsertAfter("ditor ")
[void]$body12.InsertAfter("c")
[void]$body12.InsertAfter("omments")
[void]$body12.InsertAfter("`rReplies ")
[void]$body12.InsertAfter("for comments")

# ---------------------------------------------------------------------
# Slide 13 ("Assignments")
# ---------------------------------------------------------------------
$slide13 = $p.Slides.Item(13)
$body13 = $slide13.Shapes.Item("Content Placeholder 2").TextFrame.TextRange

$body13.Text = "Comment boxes"
[void]$body13.InsertAfter("`rEmail notifications")
[void]$body13.InsertAfter("`rGradebook")
